# Apply updated "想去人数" (F column) counts and one status text change (G7)
# to both the "展览" and "全部类型" worksheets, matching the data refresh
# captured in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet names that hold the full event listing (they mirror each other).
$targetSheetNames = @("展览", "全部类型")

# F-column (想去人数) updates: row number -> new value
$fUpdates = @{
    6  = 207
    7  = 169
    10 = 300
    12 = 341
    13 = 1822
    20 = 4362
    22 = 316
    23 = 1179
    26 = 732
    28 = 381
    30 = 191
}

foreach ($sheetName in $targetSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    # G7: ticket status changed from sold out (已售罄) to sales stopped (已停售)
    $ws.Range("G7").Value = "已停售"
}
